$wb = $excel.ActiveWorkbook

# Worksheet "compounds" - update the CIViC data-source version string
$wsCompounds = $wb.Worksheets.Item("compounds")
$wsCompounds.Range("E3").Value = "24.12e"

# Worksheet "biomarkers" unaffected in content (value remains "24Q4")
$wsBiomarkers = $wb.Worksheets.Item("biomarkers")

# Update the active sheet / selection: "compounds" becomes the active tab,
# with the active cell/selection moved to E3.
$wsCompounds.Activate()
$wsCompounds.Range("E3").Select()
